$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 1884.625
$ws.Range("I38").Value = 696.3333
$ws.Range("J38").Value = 5449.5
$ws.Range("K38").Value = 2088.9999
$ws.Range("L38").Value = 16348.5
$ws.Range("M38").Value = -1716.9999
$ws.Range("N38").Value = -17092.5

# Row 113
$ws.Range("H113").Value = 3183.1765
$ws.Range("I113").Value = 3130.2307
$ws.Range("K113").Value = 3130.2307
$ws.Range("M113").Value = 123.7692999999999

# Row 129
$ws.Range("H129").Value = 4431.778
$ws.Range("I129").Value = 1962.3334
$ws.Range("K129").Value = 5887.0002
$ws.Range("M129").Value = -887.0002000000004

# Row 132
$ws.Range("H132").Value = 1578.8918
$ws.Range("I132").Value = 1454.6
$ws.Range("K132").Value = 4363.799999999999
$ws.Range("M132").Value = -1833.799999999999

# Row 137
$ws.Range("H137").Value = 3907.9524
$ws.Range("I137").Value = 2138
$ws.Range("J137").Value = 5517
$ws.Range("K137").Value = 6414
$ws.Range("L137").Value = 16551
$ws.Range("M137").Value = -3864
$ws.Range("N137").Value = -21651

# Row 138
$ws.Range("H138").Value = 13897.25
$ws.Range("I138").Value = 10958
$ws.Range("J138").Value = 15996.714
$ws.Range("K138").Value = 32874
$ws.Range("L138").Value = 47990.142
$ws.Range("M138").Value = -27734
$ws.Range("N138").Value = -58270.142

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13788.233
$ws.Range("I32").Value = 5664.7676
$ws.Range("J32").Value = 24062.03
$ws.Range("K32").Value = 5664.7676
$ws.Range("L32").Value = 24062.03
$ws.Range("M32").Value = -5377.7676
$ws.Range("N32").Value = -24636.03

# Row 61
$ws.Range("H61").Value = 1535.5
$ws.Range("I61").Value = 1547.8667
$ws.Range("J61").Value = 1350
$ws.Range("K61").Value = 1547.8667
$ws.Range("L61").Value = 1350
$ws.Range("M61").Value = -1335.8667
$ws.Range("N61").Value = -1774

# Row 74
$ws.Range("H74").Value = 4386.8335
$ws.Range("I74").Value = 1995.6
$ws.Range("J74").Value = 6094.857
$ws.Range("K74").Value = 1995.6
$ws.Range("L74").Value = 6094.857
$ws.Range("M74").Value = -1121.6
$ws.Range("N74").Value = -7842.857

# Row 77
$ws.Range("H77").Value = 4386.8335
$ws.Range("I77").Value = 1995.6
$ws.Range("J77").Value = 6094.857
$ws.Range("K77").Value = 9978
$ws.Range("L77").Value = 30474.285
$ws.Range("M77").Value = -5610
$ws.Range("N77").Value = -39210.285

# Row 97
$ws.Range("H97").Value = 398.66666
$ws.Range("I97").Value = 437.3
$ws.Range("J97").Value = 205.5
$ws.Range("K97").Value = 437.3
$ws.Range("L97").Value = 205.5
$ws.Range("M97").Value = 58.69999999999999
$ws.Range("N97").Value = -1197.5

# Row 110
$ws.Range("H110").Value = 6717.7144
$ws.Range("I110").Value = 6717.7144
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 6717.7144
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -4672.7144

# Row 136
$ws.Range("H136").Value = 1535.5
$ws.Range("I136").Value = 1547.8667
$ws.Range("J136").Value = 1350
$ws.Range("K136").Value = 4643.6001
$ws.Range("L136").Value = 4050
$ws.Range("M136").Value = -2093.6001
$ws.Range("N136").Value = -9150

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1787.7407
$ws.Range("I107").Value = 1322.25
$ws.Range("K107").Value = 1322.25
$ws.Range("M107").Value = 597.75

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4005
$ws.Range("J16").Value = 4999
$ws.Range("L16").Value = 4999
$ws.Range("N16").Value = -5573

# Row 31
$ws.Range("H31").Value = 5457.4116
$ws.Range("I31").Value = 4878
$ws.Range("K31").Value = 4878
$ws.Range("M31").Value = -4583

# Row 34
$ws.Range("H34").Value = 5457.4116
$ws.Range("I34").Value = 4878
$ws.Range("K34").Value = 4878
$ws.Range("M34").Value = -4676

# Row 62
$ws.Range("H62").Value = 48823.89
$ws.Range("J62").Value = 136664.33
$ws.Range("L62").Value = 136664.33
$ws.Range("N62").Value = -137912.33

# Row 65
$ws.Range("H65").Value = 48823.89
$ws.Range("J65").Value = 136664.33
$ws.Range("L65").Value = 683321.6499999999
$ws.Range("N65").Value = -689561.6499999999

# Row 86
$ws.Range("H86").Value = 15602.6
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877

# Row 89
$ws.Range("H89").Value = 15602.6
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384

# Row 103
$ws.Range("H103").Value = 13056.091
$ws.Range("I103").Value = 11361.7
$ws.Range("K103").Value = 11361.7
$ws.Range("M103").Value = -10189.7

# Row 113
$ws.Range("H113").Value = 4005
$ws.Range("J113").Value = 4999
$ws.Range("L113").Value = 4999
$ws.Range("N113").Value = -9339

# Row 132
$ws.Range("H132").Value = 3567.25
$ws.Range("J132").Value = 5541
$ws.Range("L132").Value = 16623
$ws.Range("N132").Value = -21683

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 282
$ws.Range("J23").Value = 250
$ws.Range("L23").Value = 750
$ws.Range("N23").Value = -1220

# Row 38
$ws.Range("H38").Value = 180.86363
$ws.Range("I38").Value = 98.38461
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 295.15383
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = 51.84617000000003
$ws.Range("N38").Value = -1594

# Row 131
$ws.Range("H131").Value = 1499.5
$ws.Range("I131").Value = 749.75
$ws.Range("J131").Value = 2999
$ws.Range("K131").Value = 2249.25
$ws.Range("L131").Value = 8997
$ws.Range("M131").Value = 2790.75
$ws.Range("N131").Value = -19077

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 426541.03
$ws.Range("I122").Value = 113521.89
$ws.Range("K122").Value = 340565.67
$ws.Range("M122").Value = -338115.67

# Row 126
$ws.Range("H126").Value = 4316.5386
$ws.Range("I126").Value = 3454.4
$ws.Range("J126").Value = 4855.375
$ws.Range("K126").Value = 10363.2
$ws.Range("L126").Value = 14566.125
$ws.Range("M126").Value = -7893.200000000001
$ws.Range("N126").Value = -19506.125

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3925.3333
$ws.Range("I7").Value = 3981.4
$ws.Range("K7").Value = 3981.4
$ws.Range("M7").Value = -3869.4

# Row 16
$ws.Range("H16").Value = 8517.75
$ws.Range("I16").Value = 8517.75
$ws.Range("K16").Value = 8517.75
$ws.Range("M16").Value = -8347.75

# Row 61
$ws.Range("H61").Value = 5027.6
$ws.Range("I61").Value = 4547.3335
$ws.Range("K61").Value = 4547.3335
$ws.Range("M61").Value = -4345.3335

# Row 113
$ws.Range("H113").Value = 5027.6
$ws.Range("I113").Value = 4547.3335
$ws.Range("K113").Value = 4547.3335
$ws.Range("M113").Value = -2377.3335

# Row 126
$ws.Range("H126").Value = 3925.3333
$ws.Range("I126").Value = 3981.4
$ws.Range("K126").Value = 11944.2
$ws.Range("M126").Value = -9474.200000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Row 132
$ws.Range("H132").Value = 1390.8966
$ws.Range("I132").Value = 1053.52
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 3160.56
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -630.5599999999999
$ws.Range("N132").Value = -15558.5

# Row 136
$ws.Range("H136").Value = 55848.844
$ws.Range("I136").Value = 2809
$ws.Range("K136").Value = 8427
$ws.Range("M136").Value = -5877
